{"js": "// Apply the Maze Escape Megaways copy edits.\n// Each entry: [searchText, replacementText]\nconst replacements = [\n  [\n    \"Play Free Maze Escape Megaways Slot Game Review\",\n    \"Play Maze Escape Megaways - Free Slot Game\"\n  ],\n  [\n    \"Innovative gameplay mechanics and features\",\n    \"Innovative gameplay mechanics with cascading reels and shifting maze avalanches\"\n  ],\n  [\n    \"Attractive maximum prize of 25,000x the total stake\",\n    \"Special symbols unlock stacked wilds and multipliers up to 7x\"\n  ],\n  [\n    \"Greek mythology theme with immersive setting\",\n    \"Colorful, cartoon-style graphics with an Ancient Greece-themed setting\"\n  ],\n  [\n    \"High RTP of 96.15%, slightly above average\",\n    \"High volatility and above-average RTP of 96.15%\"\n  ],\n  [\n    \"High volatility may not appeal to all players\",\n    \"Limited number of bonus symbols and special features\"\n  ],\n  [\n    \"Limited bonus symbols and special features\",\n    \"May not appeal to players who prefer traditional slot game formats\"\n  ],\n  [\n    \"Read our review of Maze Escape Megaways slot game, play for free and experience an innovative Megaways gameplay with high RTP and attractive max prize.\",\n    \"Experience the innovative gameplay and Greek mythology theme in Maze Escape Megaways. Play for free!\"\n  ]\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Maze Escape Megaways copy edits.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Free Maze Escape Megaways Slot Game Review\"; Replace = \"Play Maze Escape Megaways - Free Slot Game\" },\n    @{ Find = \"Innovative gameplay mechanics and features\"; Replace = \"Innovative gameplay mechanics with cascading reels and shifting maze avalanches\" },\n    @{ Find = \"Attractive maximum prize of 25,000x the total stake\"; Replace = \"Special symbols unlock stacked wilds and multipliers up to 7x\" },\n    @{ Find = \"Greek mythology theme with immersive setting\"; Replace = \"Colorful, cartoon-style graphics with an Ancient Greece-themed setting\" },\n    @{ Find = \"High RTP of 96.15%, slightly above average\"; Replace = \"High volatility and above-average RTP of 96.15%\" },\n    @{ Find = \"High volatility may not appeal to all players\"; Replace = \"Limited number of bonus symbols and special features\" },\n    @{ Find = \"Limited bonus symbols and special features\"; Replace = \"May not appeal to players who prefer traditional slot game formats\" },\n    @{ Find = \"Read our review of Maze Escape Megaways slot game, play for free and experience an innovative Megaways gameplay with high RTP and attractive max prize.\"; Replace = \"Experience the innovative gameplay and Greek mythology theme in Maze Escape Megaways. Play for free!\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$r.Find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$r.Replace, [ref]2) | Out-Null\n}\n"}
